$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (Test Case 4): result moves from "OK" to a new intermediate "Medio OK" state
$ws.Range("F5").Value = "Medio OK"
$ws.Range("F5").Interior.Color = 13431551

# Row 6 (Test Case 5): result moves from "OK" back to "PTE" (reuse the PTE formatting)
$ws.Range("F6").Value = "PTE"
$ws.Range("F7").Copy()
$ws.Range("F6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Rows 9, 10, 11, 13 (Test Cases 8, 9, 10, 12): the sales loaded from the database
# are now shown, so these test cases move from "PTE" to "OK" (reuse the OK formatting)
$ws.Range("F2").Copy()

$ws.Range("F9").Value = "OK"
$ws.Range("F9").PasteSpecial(-4122)

$ws.Range("F10").Value = "OK"
$ws.Range("F10").PasteSpecial(-4122)

$ws.Range("F11").Value = "OK"
$ws.Range("F11").PasteSpecial(-4122)

$ws.Range("F13").Value = "OK"
$ws.Range("F13").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# A13 is the last row of the merged A9:A13 block; give it the same bottom border
# treatment as A8 (last row of the merged A2:A8 block) so the block is closed off
$ws.Range("A8").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Cursor was left on F5 when the file was saved
$ws.Range("F5").Select()
